# Update "Results" sheet rows 2-9 (data rows), reflecting the re-optimized
# supplier assignments, DDP freight consolidation (freight folded into unit
# price; Freight Amount column zeroed out), and new discount tiers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Row 2
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = 'A'
$ws.Cells.Item(2,3).Value = 'Brunswick'
$ws.Cells.Item(2,4).Value = 'Supplier 2'
$ws.Cells.Item(2,5).Value = 25
$ws.Cells.Item(2,6).Value = 36400
$ws.Cells.Item(2,7).Value = 'Supplier 2'
$ws.Cells.Item(2,8).Value = 14
$ws.Cells.Item(2,9).Value = '10%'
$ws.Cells.Item(2,10).Value = 12.6
$ws.Cells.Item(2,11).Value = 'DDP'
$ws.Cells.Item(2,12).Value = 0
$ws.Cells.Item(2,13).Value = 12.6
$ws.Cells.Item(2,14).Value = 18345.6
$ws.Cells.Item(2,15).Value = 1456
$ws.Cells.Item(2,16).Value = 18054.4
$ws.Cells.Item(2,17).Value = '0%'
$ws.Cells.Item(2,18).Value = 0

# Row 3
$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = 'A'
$ws.Cells.Item(3,3).Value = 'Brunswick'
$ws.Cells.Item(3,4).Value = 'Supplier 3'
$ws.Cells.Item(3,5).Value = 36
$ws.Cells.Item(3,6).Value = 85608
$ws.Cells.Item(3,7).Value = 'Supplier 4'
$ws.Cells.Item(3,8).Value = 24
$ws.Cells.Item(3,9).Value = '0%'
$ws.Cells.Item(3,10).Value = 24
$ws.Cells.Item(3,11).Value = 'DDP'
$ws.Cells.Item(3,12).Value = 0
$ws.Cells.Item(3,13).Value = 24
$ws.Cells.Item(3,14).Value = 57072
$ws.Cells.Item(3,15).Value = 2378
$ws.Cells.Item(3,16).Value = 28536
$ws.Cells.Item(3,17).Value = '0%'
$ws.Cells.Item(3,18).Value = 0

# Row 4
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 'B'
$ws.Cells.Item(4,3).Value = 'Brunswick'
$ws.Cells.Item(4,4).Value = 'Supplier 3'
$ws.Cells.Item(4,5).Value = 36
$ws.Cells.Item(4,6).Value = 1548
$ws.Cells.Item(4,7).Value = 'Supplier 2'
$ws.Cells.Item(4,8).Value = 15
$ws.Cells.Item(4,9).Value = '10%'
$ws.Cells.Item(4,10).Value = 13.5
$ws.Cells.Item(4,11).Value = 'DDP'
$ws.Cells.Item(4,12).Value = 0
$ws.Cells.Item(4,13).Value = 13.5
$ws.Cells.Item(4,14).Value = 580.5
$ws.Cells.Item(4,15).Value = 43
$ws.Cells.Item(4,16).Value = 967.5
$ws.Cells.Item(4,17).Value = '0%'
$ws.Cells.Item(4,18).Value = 0

# Row 5
$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 'A'
$ws.Cells.Item(5,3).Value = 'Palatka'
$ws.Cells.Item(5,4).Value = 'Supplier 2'
$ws.Cells.Item(5,5).Value = 82
$ws.Cells.Item(5,6).Value = 277324
$ws.Cells.Item(5,7).Value = 'Supplier 3'
$ws.Cells.Item(5,8).Value = 56
$ws.Cells.Item(5,9).Value = '0%'
$ws.Cells.Item(5,10).Value = 56
$ws.Cells.Item(5,11).Value = 'DDP'
$ws.Cells.Item(5,12).Value = 0
$ws.Cells.Item(5,13).Value = 56
$ws.Cells.Item(5,14).Value = 189392
$ws.Cells.Item(5,15).Value = 3382
$ws.Cells.Item(5,16).Value = 87932
$ws.Cells.Item(5,17).Value = '0%'
$ws.Cells.Item(5,18).Value = 0

# Row 6
$ws.Cells.Item(6,1).Value = 3
$ws.Cells.Item(6,2).Value = 'B'
$ws.Cells.Item(6,3).Value = 'Palatka'
$ws.Cells.Item(6,4).Value = 'Supplier 2'
$ws.Cells.Item(6,5).Value = 82
$ws.Cells.Item(6,6).Value = 252560
$ws.Cells.Item(6,7).Value = 'Supplier 1'
$ws.Cells.Item(6,8).Value = 15
$ws.Cells.Item(6,9).Value = '5%'
$ws.Cells.Item(6,10).Value = 14.25
$ws.Cells.Item(6,11).Value = 'DDP'
$ws.Cells.Item(6,12).Value = 0
$ws.Cells.Item(6,13).Value = 14.25
$ws.Cells.Item(6,14).Value = 43890
$ws.Cells.Item(6,15).Value = 3080
$ws.Cells.Item(6,16).Value = 208670
$ws.Cells.Item(6,17).Value = '0%'
$ws.Cells.Item(6,18).Value = 0

# Row 7
$ws.Cells.Item(7,1).Value = 4
$ws.Cells.Item(7,2).Value = 'A'
$ws.Cells.Item(7,3).Value = 'Big Island'
$ws.Cells.Item(7,4).Value = 'Supplier 1'
$ws.Cells.Item(7,5).Value = 32
$ws.Cells.Item(7,6).Value = 75424
$ws.Cells.Item(7,7).Value = 'Supplier 5'
$ws.Cells.Item(7,8).Value = 13
$ws.Cells.Item(7,9).Value = '0%'
$ws.Cells.Item(7,10).Value = 13
$ws.Cells.Item(7,11).Value = 'DDP'
$ws.Cells.Item(7,12).Value = 0
$ws.Cells.Item(7,13).Value = 13
$ws.Cells.Item(7,14).Value = 30641
$ws.Cells.Item(7,15).Value = 2357
$ws.Cells.Item(7,16).Value = 44783
$ws.Cells.Item(7,17).Value = '0%'
$ws.Cells.Item(7,18).Value = 0

# Row 8
$ws.Cells.Item(8,1).Value = 5
$ws.Cells.Item(8,2).Value = 'A'
$ws.Cells.Item(8,3).Value = 'Huntsville'
$ws.Cells.Item(8,4).Value = 'Supplier 3'
$ws.Cells.Item(8,5).Value = 75
$ws.Cells.Item(8,6).Value = 144000
$ws.Cells.Item(8,7).Value = 'Supplier 1'
$ws.Cells.Item(8,8).Value = 24
$ws.Cells.Item(8,9).Value = '5%'
$ws.Cells.Item(8,10).Value = 22.8
$ws.Cells.Item(8,11).Value = 'DDP'
$ws.Cells.Item(8,12).Value = 0
$ws.Cells.Item(8,13).Value = 22.8
$ws.Cells.Item(8,14).Value = 43775.99999999999
$ws.Cells.Item(8,15).Value = 1920
$ws.Cells.Item(8,16).Value = 100224
$ws.Cells.Item(8,17).Value = '0%'
$ws.Cells.Item(8,18).Value = 0

# Row 9
$ws.Cells.Item(9,1).Value = 5
$ws.Cells.Item(9,2).Value = 'B'
$ws.Cells.Item(9,3).Value = 'Huntsville'
$ws.Cells.Item(9,4).Value = 'Supplier 3'
$ws.Cells.Item(9,5).Value = 75
$ws.Cells.Item(9,6).Value = 37575
$ws.Cells.Item(9,7).Value = 'Supplier 2'
$ws.Cells.Item(9,8).Value = 15
$ws.Cells.Item(9,9).Value = '10%'
$ws.Cells.Item(9,10).Value = 13.5
$ws.Cells.Item(9,11).Value = 'DDP'
$ws.Cells.Item(9,12).Value = 0
$ws.Cells.Item(9,13).Value = 13.5
$ws.Cells.Item(9,14).Value = 6763.5
$ws.Cells.Item(9,15).Value = 501
$ws.Cells.Item(9,16).Value = 30811.5
$ws.Cells.Item(9,17).Value = '0%'
$ws.Cells.Item(9,18).Value = 0

# Update "LP Model" sheet with the regenerated LP formulation text
# (adds discount-tier constraints/binaries, drops the Fix_d_Supplier_1/2
# fixes, and zeroes out the per-unit Freight_Supplier_* equations).
$lpWs = $wb.Worksheets.Item("LP Model")
$lpText = @'
\* Sourcing_with_MultiTier_Rebates_Discounts *\
Minimize
OBJ: S_Supplier_1 + S_Supplier_2 + S_Supplier_3 + S_Supplier_4 + S_Supplier_5
 - rebate_Supplier_1 - rebate_Supplier_2 - rebate_Supplier_3
 - rebate_Supplier_4 - rebate_Supplier_5
Subject To
ActiveLink_Supplier_1_1: x_Supplier_1_1 - 1000000000 z_Supplier_1 <= 0
ActiveLink_Supplier_1_2: x_Supplier_1_2 - 1000000000 z_Supplier_1 <= 0
ActiveLink_Supplier_1_3: x_Supplier_1_3 - 1000000000 z_Supplier_1 <= 0
ActiveLink_Supplier_1_4: x_Supplier_1_4 - 1000000000 z_Supplier_1 <= 0
ActiveLink_Supplier_1_5: x_Supplier_1_5 - 1000000000 z_Supplier_1 <= 0
ActiveLink_Supplier_2_1: x_Supplier_2_1 - 1000000000 z_Supplier_2 <= 0
ActiveLink_Supplier_2_2: x_Supplier_2_2 - 1000000000 z_Supplier_2 <= 0
ActiveLink_Supplier_2_3: x_Supplier_2_3 - 1000000000 z_Supplier_2 <= 0
ActiveLink_Supplier_2_4: x_Supplier_2_4 - 1000000000 z_Supplier_2 <= 0
ActiveLink_Supplier_2_5: x_Supplier_2_5 - 1000000000 z_Supplier_2 <= 0
ActiveLink_Supplier_3_1: x_Supplier_3_1 - 1000000000 z_Supplier_3 <= 0
ActiveLink_Supplier_3_2: x_Supplier_3_2 - 1000000000 z_Supplier_3 <= 0
ActiveLink_Supplier_3_3: x_Supplier_3_3 - 1000000000 z_Supplier_3 <= 0
ActiveLink_Supplier_3_4: x_Supplier_3_4 - 1000000000 z_Supplier_3 <= 0
ActiveLink_Supplier_3_5: x_Supplier_3_5 - 1000000000 z_Supplier_3 <= 0
ActiveLink_Supplier_4_1: x_Supplier_4_1 - 1000000000 z_Supplier_4 <= 0
ActiveLink_Supplier_4_2: x_Supplier_4_2 - 1000000000 z_Supplier_4 <= 0
ActiveLink_Supplier_4_3: x_Supplier_4_3 - 1000000000 z_Supplier_4 <= 0
ActiveLink_Supplier_4_4: x_Supplier_4_4 - 1000000000 z_Supplier_4 <= 0
ActiveLink_Supplier_4_5: x_Supplier_4_5 - 1000000000 z_Supplier_4 <= 0
ActiveLink_Supplier_5_1: x_Supplier_5_1 - 1000000000 z_Supplier_5 <= 0
ActiveLink_Supplier_5_2: x_Supplier_5_2 - 1000000000 z_Supplier_5 <= 0
ActiveLink_Supplier_5_3: x_Supplier_5_3 - 1000000000 z_Supplier_5 <= 0
ActiveLink_Supplier_5_4: x_Supplier_5_4 - 1000000000 z_Supplier_5 <= 0
ActiveLink_Supplier_5_5: x_Supplier_5_5 - 1000000000 z_Supplier_5 <= 0
BaseSpend_Supplier_1: S0_Supplier_1 - 34 x_Supplier_1_1 - 15 x_Supplier_1_2
 - 15 x_Supplier_1_3 - 75 x_Supplier_1_4 - 24 x_Supplier_1_5 = 0
BaseSpend_Supplier_2: S0_Supplier_2 - 14 x_Supplier_2_1 - 15 x_Supplier_2_2
 - 78 x_Supplier_2_3 - 34 x_Supplier_2_4 - 15 x_Supplier_2_5 = 0
BaseSpend_Supplier_3: S0_Supplier_3 - 75 x_Supplier_3_1 - 25 x_Supplier_3_2
 - 56 x_Supplier_3_3 - 24 x_Supplier_3_4 - 87 x_Supplier_3_5 = 0
BaseSpend_Supplier_4: S0_Supplier_4 - 93 x_Supplier_4_1 - 24 x_Supplier_4_2
 - 78 x_Supplier_4_3 - 56 x_Supplier_4_4 - 86 x_Supplier_4_5 = 0
BaseSpend_Supplier_5: S0_Supplier_5 - 76 x_Supplier_5_1 - 32 x_Supplier_5_2
 - 89 x_Supplier_5_3 - 13 x_Supplier_5_4 - 68 x_Supplier_5_5 = 0
Capacity_Supplier_1_Bid_ID_1: x_Supplier_1_1 <= 1457
Capacity_Supplier_1_Bid_ID_2: x_Supplier_1_2 <= 2422
Capacity_Supplier_1_Bid_ID_3: x_Supplier_1_3 <= 6463
Capacity_Supplier_1_Bid_ID_4: x_Supplier_1_4 <= 2358
Capacity_Supplier_1_Bid_ID_5: x_Supplier_1_5 <= 2422
Capacity_Supplier_2_Capacity_Group_Category_1: x_Supplier_2_1 + x_Supplier_2_3
 + x_Supplier_2_4 <= 8001
Capacity_Supplier_2_Capacity_Group_Category_2: x_Supplier_2_2 <= 2501
Capacity_Supplier_2_Capacity_Group_Category_3: x_Supplier_2_5 <= 501
Capacity_Supplier_3_Bid_ID_1: x_Supplier_3_1 <= 1457
Capacity_Supplier_3_Bid_ID_2: x_Supplier_3_2 <= 2422
Capacity_Supplier_3_Bid_ID_3: x_Supplier_3_3 <= 6463
Capacity_Supplier_3_Bid_ID_4: x_Supplier_3_4 <= 2358
Capacity_Supplier_3_Bid_ID_5: x_Supplier_3_5 <= 2422
Capacity_Supplier_4_Bid_ID_1: x_Supplier_4_1 <= 1457
Capacity_Supplier_4_Bid_ID_2: x_Supplier_4_2 <= 2422
Capacity_Supplier_4_Bid_ID_3: x_Supplier_4_3 <= 6463
Capacity_Supplier_4_Bid_ID_4: x_Supplier_4_4 <= 2358
Capacity_Supplier_4_Bid_ID_5: x_Supplier_4_5 <= 2422
Capacity_Supplier_5_Description_Large_Item: x_Supplier_5_5 <= 5001
Capacity_Supplier_5_Description_Medium_item: x_Supplier_5_2 <= 2501
Capacity_Supplier_5_Description_Small_item: x_Supplier_5_1 + x_Supplier_5_3
 + x_Supplier_5_4 <= 9001
Demand_1: x_Supplier_1_1 + x_Supplier_2_1 + x_Supplier_3_1 + x_Supplier_4_1
 + x_Supplier_5_1 = 1456
Demand_2: x_Supplier_1_2 + x_Supplier_2_2 + x_Supplier_3_2 + x_Supplier_4_2
 + x_Supplier_5_2 = 2421
Demand_3: x_Supplier_1_3 + x_Supplier_2_3 + x_Supplier_3_3 + x_Supplier_4_3
 + x_Supplier_5_3 = 6462
Demand_4: x_Supplier_1_4 + x_Supplier_2_4 + x_Supplier_3_4 + x_Supplier_4_4
 + x_Supplier_5_4 = 2357
Demand_5: x_Supplier_1_5 + x_Supplier_2_5 + x_Supplier_3_5 + x_Supplier_4_5
 + x_Supplier_5_5 = 2421
DiscountTierLower_Supplier_1_0: - 0.05 S0_Supplier_1 + d_Supplier_1
 - 1406346 z_discount_Supplier_1_0 >= -1406346
DiscountTierLower_Supplier_2_0: - 0.1 S0_Supplier_2 + d_Supplier_2
 - 1023279 z_discount_Supplier_2_0 >= -1023279
DiscountTierMax_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3
 + x_Supplier_1_4 + x_Supplier_1_5 + 1406346 z_discount_Supplier_1_0
 <= 1411346
DiscountTierMax_Supplier_2_0: x_Supplier_2_1 + x_Supplier_2_2 + x_Supplier_2_3
 + x_Supplier_2_4 + x_Supplier_2_5 + 1023279 z_discount_Supplier_2_0
 <= 1025279
DiscountTierMin_Supplier_1_0: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3
 + x_Supplier_1_4 + x_Supplier_1_5 >= 0
DiscountTierMin_Supplier_2_0: x_Supplier_2_1 + x_Supplier_2_2 + x_Supplier_2_3
 + x_Supplier_2_4 + x_Supplier_2_5 >= 0
DiscountTierSelect_Supplier_1: z_discount_Supplier_1_0 = 1
DiscountTierSelect_Supplier_2: z_discount_Supplier_2_0 = 1
DiscountTierUpper_Supplier_1_0: - 0.05 S0_Supplier_1 + d_Supplier_1
 + 1406346 z_discount_Supplier_1_0 <= 1406346
DiscountTierUpper_Supplier_2_0: - 0.1 S0_Supplier_2 + d_Supplier_2
 + 1023279 z_discount_Supplier_2_0 <= 1023279
EffectiveSpend_Supplier_1: - F_Supplier_1 - S0_Supplier_1 + S_Supplier_1
 + d_Supplier_1 = 0
EffectiveSpend_Supplier_2: - F_Supplier_2 - S0_Supplier_2 + S_Supplier_2
 + d_Supplier_2 = 0
EffectiveSpend_Supplier_3: - F_Supplier_3 - S0_Supplier_3 + S_Supplier_3
 + d_Supplier_3 = 0
EffectiveSpend_Supplier_4: - F_Supplier_4 - S0_Supplier_4 + S_Supplier_4
 + d_Supplier_4 = 0
EffectiveSpend_Supplier_5: - F_Supplier_5 - S0_Supplier_5 + S_Supplier_5
 + d_Supplier_5 = 0
Fix_d_Supplier_3: d_Supplier_3 = 0
Fix_d_Supplier_4: d_Supplier_4 = 0
Fix_d_Supplier_5: d_Supplier_5 = 0
Fix_rebate_Supplier_1: rebate_Supplier_1 = 0
Fix_rebate_Supplier_2: rebate_Supplier_2 = 0
Fix_rebate_Supplier_3: rebate_Supplier_3 = 0
Fix_rebate_Supplier_4: rebate_Supplier_4 = 0
Fix_rebate_Supplier_5: rebate_Supplier_5 = 0
Freight_Supplier_1: F_Supplier_1 = 0
Freight_Supplier_2: F_Supplier_2 = 0
Freight_Supplier_3: F_Supplier_3 = 0
Freight_Supplier_4: F_Supplier_4 = 0
Freight_Supplier_5: F_Supplier_5 = 0
MinAward_Supplier_1: x_Supplier_1_1 + x_Supplier_1_2 + x_Supplier_1_3
 + x_Supplier_1_4 + x_Supplier_1_5 - z_Supplier_1 >= 0
MinAward_Supplier_2: x_Supplier_2_1 + x_Supplier_2_2 + x_Supplier_2_3
 + x_Supplier_2_4 + x_Supplier_2_5 - z_Supplier_2 >= 0
MinAward_Supplier_3: x_Supplier_3_1 + x_Supplier_3_2 + x_Supplier_3_3
 + x_Supplier_3_4 + x_Supplier_3_5 - z_Supplier_3 >= 0
MinAward_Supplier_4: x_Supplier_4_1 + x_Supplier_4_2 + x_Supplier_4_3
 + x_Supplier_4_4 + x_Supplier_4_5 - z_Supplier_4 >= 0
MinAward_Supplier_5: x_Supplier_5_1 + x_Supplier_5_2 + x_Supplier_5_3
 + x_Supplier_5_4 + x_Supplier_5_5 - z_Supplier_5 >= 0
TransitionLower_1_Supplier_1: x_Supplier_1_1 >= 0
TransitionLower_1_Supplier_3: x_Supplier_3_1 >= 0
TransitionLower_1_Supplier_4: x_Supplier_4_1 >= 0
TransitionLower_1_Supplier_5: x_Supplier_5_1 >= 0
TransitionLower_2_Supplier_1: x_Supplier_1_2 >= 0
TransitionLower_2_Supplier_2: x_Supplier_2_2 >= 0
TransitionLower_2_Supplier_4: x_Supplier_4_2 >= 0
TransitionLower_2_Supplier_5: x_Supplier_5_2 >= 0
TransitionLower_3_Supplier_1: x_Supplier_1_3 >= 0
TransitionLower_3_Supplier_3: x_Supplier_3_3 >= 0
TransitionLower_3_Supplier_4: x_Supplier_4_3 >= 0
TransitionLower_3_Supplier_5: x_Supplier_5_3 >= 0
TransitionLower_4_Supplier_2: x_Supplier_2_4 >= 0
TransitionLower_4_Supplier_3: x_Supplier_3_4 >= 0
TransitionLower_4_Supplier_4: x_Supplier_4_4 >= 0
TransitionLower_4_Supplier_5: x_Supplier_5_4 >= 0
TransitionLower_5_Supplier_1: x_Supplier_1_5 >= 0
TransitionLower_5_Supplier_2: x_Supplier_2_5 >= 0
TransitionLower_5_Supplier_4: x_Supplier_4_5 >= 0
TransitionLower_5_Supplier_5: x_Supplier_5_5 >= 0
Transition_1_Supplier_1: - 1456 T_1_Supplier_1 + x_Supplier_1_1 <= 0
Transition_1_Supplier_3: - 1456 T_1_Supplier_3 + x_Supplier_3_1 <= 0
Transition_1_Supplier_4: - 1456 T_1_Supplier_4 + x_Supplier_4_1 <= 0
Transition_1_Supplier_5: - 1456 T_1_Supplier_5 + x_Supplier_5_1 <= 0
Transition_2_Supplier_1: - 2421 T_2_Supplier_1 + x_Supplier_1_2 <= 0
Transition_2_Supplier_2: - 2421 T_2_Supplier_2 + x_Supplier_2_2 <= 0
Transition_2_Supplier_4: - 2421 T_2_Supplier_4 + x_Supplier_4_2 <= 0
Transition_2_Supplier_5: - 2421 T_2_Supplier_5 + x_Supplier_5_2 <= 0
Transition_3_Supplier_1: - 6462 T_3_Supplier_1 + x_Supplier_1_3 <= 0
Transition_3_Supplier_3: - 6462 T_3_Supplier_3 + x_Supplier_3_3 <= 0
Transition_3_Supplier_4: - 6462 T_3_Supplier_4 + x_Supplier_4_3 <= 0
Transition_3_Supplier_5: - 6462 T_3_Supplier_5 + x_Supplier_5_3 <= 0
Transition_4_Supplier_2: - 2357 T_4_Supplier_2 + x_Supplier_2_4 <= 0
Transition_4_Supplier_3: - 2357 T_4_Supplier_3 + x_Supplier_3_4 <= 0
Transition_4_Supplier_4: - 2357 T_4_Supplier_4 + x_Supplier_4_4 <= 0
Transition_4_Supplier_5: - 2357 T_4_Supplier_5 + x_Supplier_5_4 <= 0
Transition_5_Supplier_1: - 2421 T_5_Supplier_1 + x_Supplier_1_5 <= 0
Transition_5_Supplier_2: - 2421 T_5_Supplier_2 + x_Supplier_2_5 <= 0
Transition_5_Supplier_4: - 2421 T_5_Supplier_4 + x_Supplier_4_5 <= 0
Transition_5_Supplier_5: - 2421 T_5_Supplier_5 + x_Supplier_5_5 <= 0
Volume_Supplier_1: V_Supplier_1 - x_Supplier_1_1 - x_Supplier_1_2
 - x_Supplier_1_3 - x_Supplier_1_4 - x_Supplier_1_5 = 0
Volume_Supplier_2: V_Supplier_2 - x_Supplier_2_1 - x_Supplier_2_2
 - x_Supplier_2_3 - x_Supplier_2_4 - x_Supplier_2_5 = 0
Volume_Supplier_3: V_Supplier_3 - x_Supplier_3_1 - x_Supplier_3_2
 - x_Supplier_3_3 - x_Supplier_3_4 - x_Supplier_3_5 = 0
Volume_Supplier_4: V_Supplier_4 - x_Supplier_4_1 - x_Supplier_4_2
 - x_Supplier_4_3 - x_Supplier_4_4 - x_Supplier_4_5 = 0
Volume_Supplier_5: V_Supplier_5 - x_Supplier_5_1 - x_Supplier_5_2
 - x_Supplier_5_3 - x_Supplier_5_4 - x_Supplier_5_5 = 0
Binaries
T_1_Supplier_1
T_1_Supplier_3
T_1_Supplier_4
T_1_Supplier_5
T_2_Supplier_1
T_2_Supplier_2
T_2_Supplier_4
T_2_Supplier_5
T_3_Supplier_1
T_3_Supplier_3
T_3_Supplier_4
T_3_Supplier_5
T_4_Supplier_2
T_4_Supplier_3
T_4_Supplier_4
T_4_Supplier_5
T_5_Supplier_1
T_5_Supplier_2
T_5_Supplier_4
T_5_Supplier_5
z_Supplier_1
z_Supplier_2
z_Supplier_3
z_Supplier_4
z_Supplier_5
z_discount_Supplier_1_0
z_discount_Supplier_2_0
End

'@
$lpWs.Cells.Item(2,1).Value = $lpText
